$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Pre-register a bold Calibri 11 font in the style table (matching the
#     font used later for the partial rich-text run) and a wrap-text cell
#     style, without leaving stray formatting on the sheet. A scratch cell
#     far outside the used range is formatted then the whole row is removed.
$scratchRow = 50
$scratch = $ws.Cells.Item($scratchRow,1)
$scratch.Value2 = "x"
$scratch.WrapText = $true
$scratch.Font.Bold = $true
$scratch.Font.Bold = $false

# Column A width
$ws.Columns.Item(1).ColumnWidth = 99.45182291666667

$ws.Cells.Item(1,1).Value2 = "EXTERNAL PROGRAMS"

$ws.Cells.Item(2,1).Value2 = "tn:Link external programs that print temperature when called.  This allows to connect meters that use any program language.`nArtisan will start the program each sample period.  The program output must be to Stdout (like when using print statements).  The program must exit and must not be persistent."
$ws.Cells.Item(2,1).WrapText = $true
$ws.Rows.Item(2).RowHeight = 75

$ws.Cells.Item(3,1).Value2 = "tn:"

$ws.Cells.Item(4,1).Value2 = "tn:If only one termperature is provided it will be interpreted as BT.  If more than one temperature is provided the values are order dependent with ET first and BT second."
$ws.Cells.Item(4,1).WrapText = $true
$ws.Rows.Item(4).RowHeight = 30

$ws.Cells.Item(5,1).Value2 = "tn:"
$ws.Cells.Item(5,1).WrapText = $true

$ws.Cells.Item(6,1).Value2 = "tn:Data may also be provided to the `"Program`" extra devices.  Extra device `"Program`" are the first two values, typically ET and BT.  `"Program 34`" are the third and fourth values.  Up to 10 values may be supplied."

$ws.Cells.Item(7,1).Value2 = "tn:"
$ws.Cells.Item(7,1).WrapText = $true

$ws.Cells.Item(8,1).Value2 = "tn:"
$ws.Cells.Item(8,1).WrapText = $true

$ws.Cells.Item(9,1).Value2 = "tn:Example of output needed from program for single temperature (BT):`n`"100.4`" (note: `"`" not needed)"
$ws.Cells.Item(9,1).WrapText = $true
$ws.Rows.Item(9).RowHeight = 30

$ws.Cells.Item(10,1).Value2 = "tn:"

$ws.Cells.Item(11,1).Value2 = "tn:Example of output needed from program for double temperature (ET,BT)`n`"200.4,100.4`" (note: temperatures are separated by a comma `"ET,BT`")"
$ws.Cells.Item(11,1).WrapText = $true
$ws.Rows.Item(11).RowHeight = 30

$ws.Cells.Item(12,1).Value2 = "tn:"

$ws.Cells.Item(13,1).Value2 = "tn:Example of output needed from program for double temperature (ET,BT) and extra devices (Program and Program 34)`n`"200.4,100.4,312.4,345.6,299.0,275.5`""
$ws.Cells.Item(13,1).WrapText = $true
$ws.Rows.Item(13).RowHeight = 45

$ws.Cells.Item(14,1).Value2 = "bn:Example of a file written in python language called test.py:"

$ws.Cells.Item(15,1).Value2 = "bn:"
$ws.Cells.Item(15,1).WrapText = $true

$ws.Cells.Item(16,1).Value2 = "bn:#comment: print a string with two numbers separated by a comma"

$ws.Cells.Item(17,1).Value2 = "bn:#!/usr/bin/env python"

$ws.Cells.Item(18,1).Value2 = "bn:print (`"237.1,100.4`")"

$ws.Cells.Item(19,1).Value2 = "bn:"

$ws.Cells.Item(20,1).Value2 = "bn:Note: In many cases the path to the Python or other language executatable should be provided along with the external program path.  On Windows it is  advised to enclose the paths with quotation marks if there are any spaces, and use forward slashes '/' in the path.`n`"C:/Python38-64/python.exe`" `"c:/scripts/test.py`""
$ws.Cells.Item(20,1).WrapText = $true
$ws.Rows.Item(20).RowHeight = 60

# Bold rich-text run within row 2 (applies only to the trailing sentence)
$boldRun = $ws.Cells.Item(2,1).Characters(201,99)
$boldRun.Font.Bold = $true

# Remove the scratch cell/row used to register the bold font + wrap style
$scratch.ClearContents()
$ws.Rows.Item($scratchRow).Delete()

# Update selection to match the post-edit cursor position
$ws.Range("A21").Select()
